$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2 -- r="2"
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.027579103851454
$ws.Cells.Item(2, 4).Value = 1.038921105394123
$ws.Cells.Item(2, 5).Value = 1.048746567530101
$ws.Cells.Item(2, 6).Value = 1.053026461645807
$ws.Cells.Item(2, 9).Value = 1.034947924550092
$ws.Cells.Item(2, 10).Value = 1.032736193644655
$ws.Cells.Item(2, 11).Value = 1.041707763553306
$ws.Cells.Item(2, 12).Value = 1.051505552442187
$ws.Cells.Item(2, 13).Value = 1.055773567335963
$ws.Cells.Item(2, 14).Value = 1.014925814505569

# Row 3 -- r="3"
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.028410388610715
$ws.Cells.Item(3, 4).Value = 1.039567240331053
$ws.Cells.Item(3, 5).Value = 1.049593509313524
$ws.Cells.Item(3, 6).Value = 1.053833206538555
$ws.Cells.Item(3, 9).Value = 1.035073291883963
$ws.Cells.Item(3, 10).Value = 1.033208330303622
$ws.Cells.Item(3, 11).Value = 1.042164443922471
$ws.Cells.Item(3, 12).Value = 1.052164511444643
$ws.Cells.Item(3, 13).Value = 1.056393290961473
$ws.Cells.Item(3, 14).Value = 1.01508360832088

# Row 4 -- r="4"
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028948948965678
$ws.Cells.Item(4, 4).Value = 1.039985772271892
$ws.Cells.Item(4, 5).Value = 1.050142837386687
$ws.Cells.Item(4, 6).Value = 1.054356171255847
$ws.Cells.Item(4, 9).Value = 1.035153273349118
$ws.Cells.Item(4, 10).Value = 1.03351384478649
$ws.Cells.Item(4, 11).Value = 1.042459698973
$ws.Cells.Item(4, 12).Value = 1.052591554686157
$ws.Cells.Item(4, 13).Value = 1.056794584527109
$ws.Cells.Item(4, 14).Value = 1.015185673932038

# Row 5 -- r="5"
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.029175516948152
$ws.Cells.Item(5, 4).Value = 1.040161826442164
$ws.Cells.Item(5, 5).Value = 1.050374084173284
$ws.Cells.Item(5, 6).Value = 1.054576250286363
$ws.Cells.Item(5, 9).Value = 1.03518662404354
$ws.Cells.Item(5, 10).Value = 1.03364228432417
$ws.Cells.Item(5, 11).Value = 1.042583763508674
$ws.Cells.Item(5, 12).Value = 1.052771238266798
$ws.Cells.Item(5, 13).Value = 1.05696335632704
$ws.Cells.Item(5, 14).Value = 1.015228572862382

# Row 6 -- r="6"
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.029213567880536
$ws.Cells.Item(6, 4).Value = 1.040191392717998
$ws.Cells.Item(6, 5).Value = 1.05041292960405
$ws.Cells.Item(6, 6).Value = 1.054613215639933
$ws.Cells.Item(6, 9).Value = 1.035192207712731
$ws.Cells.Item(6, 10).Value = 1.03366384992423
$ws.Cells.Item(6, 11).Value = 1.04260459089092
$ws.Cells.Item(6, 12).Value = 1.052801416948396
$ws.Cells.Item(6, 13).Value = 1.056991697793085
$ws.Cells.Item(6, 14).Value = 1.015235775213698

# Row 7 -- r="7"
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.028951975762476
$ws.Cells.Item(7, 4).Value = 1.039988124312606
$ws.Cells.Item(7, 5).Value = 1.050145926104862
$ws.Cells.Item(7, 6).Value = 1.054359111082376
$ws.Cells.Item(7, 9).Value = 1.035153720058318
$ws.Cells.Item(7, 10).Value = 1.033515560998556
$ws.Cells.Item(7, 11).Value = 1.042461356969606
$ws.Cells.Item(7, 12).Value = 1.052593955021558
$ws.Cells.Item(7, 13).Value = 1.056796839399121
$ws.Cells.Item(7, 14).Value = 1.015186247187372

# Row 8 -- r="8"
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.027859902052279
$ws.Cells.Item(8, 4).Value = 1.039139377617556
$ws.Cells.Item(8, 5).Value = 1.049032525083149
$ws.Cells.Item(8, 6).Value = 1.053298907605252
$ws.Cells.Item(8, 9).Value = 1.034990528567538
$ws.Cells.Item(8, 10).Value = 1.0328957515069
$ws.Cells.Item(8, 11).Value = 1.04186215127094
$ws.Cells.Item(8, 12).Value = 1.051728114579641
$ws.Cells.Item(8, 13).Value = 1.055982944593583
$ws.Cells.Item(8, 14).Value = 1.014979149156976

# Row 9 -- r="9"
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.025940683142308
$ws.Cells.Item(9, 4).Value = 1.03764722513649
$ws.Cells.Item(9, 5).Value = 1.047080615241622
$ws.Cells.Item(9, 6).Value = 1.051438034025979
$ws.Cells.Item(9, 9).Value = 1.034694267235388
$ws.Cells.Item(9, 10).Value = 1.03180370421375
$ws.Cells.Item(9, 11).Value = 1.040804438800621
$ws.Cells.Item(9, 12).Value = 1.05020746531703
$ws.Cells.Item(9, 13).Value = 1.054551059933503
$ws.Cells.Item(9, 14).Value = 1.014613949046231

# Row 10 -- r="10"
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024664767257622
$ws.Cells.Item(10, 4).Value = 1.036654885705324
$ws.Cells.Item(10, 5).Value = 1.045786214252173
$ws.Cells.Item(10, 6).Value = 1.050202506892929
$ws.Cells.Item(10, 9).Value = 1.03449095377851
$ws.Cells.Item(10, 10).Value = 1.031075839586621
$ws.Cells.Item(10, 11).Value = 1.040098146080903
$ws.Cells.Item(10, 12).Value = 1.049197208391085
$ws.Cells.Item(10, 13).Value = 1.05359811740331
$ws.Cells.Item(10, 14).Value = 1.014370330599462

# Row 11 -- r="11"
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024113147352485
$ws.Cells.Item(11, 4).Value = 1.036225791328422
$ws.Cells.Item(11, 5).Value = 1.045227379420863
$ws.Cells.Item(11, 6).Value = 1.049668734417576
$ws.Cells.Item(11, 9).Value = 1.034401548475991
$ws.Cells.Item(11, 10).Value = 1.030760721640051
$ws.Cells.Item(11, 11).Value = 1.039792058039506
$ws.Cells.Item(11, 12).Value = 1.048760609266665
$ws.Cells.Item(11, 13).Value = 1.053185894505214
$ws.Cells.Item(11, 14).Value = 1.014264810996334

# Row 12 -- r="12"
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.023908382029026
$ws.Cells.Item(12, 4).Value = 1.036066497886005
$ws.Cells.Item(12, 5).Value = 1.045020053117745
$ws.Cells.Item(12, 6).Value = 1.049470652743477
$ws.Cells.Item(12, 9).Value = 1.034368134343814
$ws.Cells.Item(12, 10).Value = 1.030643681930316
$ws.Cells.Item(12, 11).Value = 1.039678325864935
$ws.Cells.Item(12, 12).Value = 1.048598566130302
$ws.Cells.Item(12, 13).Value = 1.053032839489664
$ws.Cells.Item(12, 14).Value = 1.014225612100118

# Row 13 -- r="13"
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.023952298937308
$ws.Cells.Item(13, 4).Value = 1.03610066270747
$ws.Cells.Item(13, 5).Value = 1.04506451397521
$ws.Cells.Item(13, 6).Value = 1.049513133521704
$ws.Cells.Item(13, 9).Value = 1.034375311055186
$ws.Cells.Item(13, 10).Value = 1.030668786911859
$ws.Cells.Item(13, 11).Value = 1.039702723474409
$ws.Cells.Item(13, 12).Value = 1.048633319058492
$ws.Cells.Item(13, 13).Value = 1.053065667435369
$ws.Cells.Item(13, 14).Value = 1.014234020579702

# Row 14 -- r="14"
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024096218707698
$ws.Cells.Item(14, 4).Value = 1.036212622209008
$ws.Cells.Item(14, 5).Value = 1.045210236655606
$ws.Cells.Item(14, 6).Value = 1.049652357134887
$ws.Cells.Item(14, 9).Value = 1.034398790631795
$ws.Cells.Item(14, 10).Value = 1.030751046907572
$ws.Cells.Item(14, 11).Value = 1.039782657658245
$ws.Cells.Item(14, 12).Value = 1.048747212084644
$ws.Cells.Item(14, 13).Value = 1.053173241636443
$ws.Cells.Item(14, 14).Value = 1.014261570886654

# Row 15 -- r="15"
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024184909819199
$ws.Cells.Item(15, 4).Value = 1.036281616309912
$ws.Cells.Item(15, 5).Value = 1.045300054375098
$ws.Cells.Item(15, 6).Value = 1.049738161998917
$ws.Cells.Item(15, 9).Value = 1.034413230027113
$ws.Cells.Item(15, 10).Value = 1.030801731253387
$ws.Cells.Item(15, 11).Value = 1.039831902827063
$ws.Cells.Item(15, 12).Value = 1.048817402510389
$ws.Cells.Item(15, 13).Value = 1.053239530040872
$ws.Cells.Item(15, 14).Value = 1.01427854499529

# Row 16 -- r="16"
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.024701394716026
$ws.Cells.Item(16, 4).Value = 1.036683375987264
$ws.Cells.Item(16, 5).Value = 1.045823337203335
$ws.Cells.Item(16, 6).Value = 1.050237957471382
$ws.Cells.Item(16, 9).Value = 1.034496858524663
$ws.Cells.Item(16, 10).Value = 1.031096754136629
$ws.Cells.Item(16, 11).Value = 1.040118454807341
$ws.Cells.Item(16, 12).Value = 1.04922620206496
$ws.Cells.Item(16, 13).Value = 1.053625484010222
$ws.Cells.Item(16, 14).Value = 1.014377332970724

# Row 17 -- r="17"
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025025603482157
$ws.Cells.Item(17, 4).Value = 1.036935549898425
$ws.Cells.Item(17, 5).Value = 1.046152021712178
$ws.Cells.Item(17, 6).Value = 1.050551793762237
$ws.Cells.Item(17, 9).Value = 1.034548950375565
$ws.Cells.Item(17, 10).Value = 1.0312818290851
$ws.Cells.Item(17, 11).Value = 1.040298133045605
$ws.Cells.Item(17, 12).Value = 1.049482859668552
$ws.Cells.Item(17, 13).Value = 1.053867693039438
$ws.Cells.Item(17, 14).Value = 1.014439292035339

# Row 18 -- r="18"
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025214791845843
$ws.Cells.Item(18, 4).Value = 1.037082695969914
$ws.Cells.Item(18, 5).Value = 1.046343896825718
$ws.Cells.Item(18, 6).Value = 1.050734966631657
$ws.Cells.Item(18, 9).Value = 1.034579202513316
$ws.Cells.Item(18, 10).Value = 1.031389785126054
$ws.Cells.Item(18, 11).Value = 1.04040291127326
$ws.Cells.Item(18, 12).Value = 1.049632645488684
$ws.Cells.Item(18, 13).Value = 1.054009008672756
$ws.Cells.Item(18, 14).Value = 1.014475428666759

# Row 19 -- r="19"
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02527931419575
$ws.Cells.Item(19, 4).Value = 1.037132878638154
$ws.Cells.Item(19, 5).Value = 1.046409348188077
$ws.Cells.Item(19, 6).Value = 1.05079744370891
$ws.Cells.Item(19, 9).Value = 1.034589495273318
$ws.Cells.Item(19, 10).Value = 1.031426596147719
$ws.Cells.Item(19, 11).Value = 1.040438633649585
$ws.Cells.Item(19, 12).Value = 1.049683732402604
$ws.Cells.Item(19, 13).Value = 1.054057200249827
$ws.Cells.Item(19, 14).Value = 1.014487749793779

# Row 20 -- r="20"
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.024990810365709
$ws.Cells.Item(20, 4).Value = 1.036908488084245
$ws.Cells.Item(20, 5).Value = 1.046116740514139
$ws.Cells.Item(20, 6).Value = 1.050518109932948
$ws.Cells.Item(20, 9).Value = 1.034543375077658
$ws.Cells.Item(20, 10).Value = 1.031261971772949
$ws.Cells.Item(20, 11).Value = 1.040278857851845
$ws.Cells.Item(20, 12).Value = 1.049455314265309
$ws.Cells.Item(20, 13).Value = 1.053841702235486
$ws.Cells.Item(20, 14).Value = 1.014432644725401

# Row 21 -- r="21"
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024053834296492
$ws.Cells.Item(21, 4).Value = 1.036179650390341
$ws.Cells.Item(21, 5).Value = 1.045167318045205
$ws.Cells.Item(21, 6).Value = 1.04961135412219
$ws.Cells.Item(21, 9).Value = 1.034391882135763
$ws.Cells.Item(21, 10).Value = 1.030726823126921
$ws.Cells.Item(21, 11).Value = 1.03975912005583
$ws.Cells.Item(21, 12).Value = 1.048713669846476
$ws.Cells.Item(21, 13).Value = 1.053141561964809
$ws.Cells.Item(21, 14).Value = 1.014253458120567

# Row 22 -- r="22"
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023465478247492
$ws.Cells.Item(22, 4).Value = 1.035721930461702
$ws.Cells.Item(22, 5).Value = 1.044571824936766
$ws.Cells.Item(22, 6).Value = 1.049042313222921
$ws.Cells.Item(22, 9).Value = 1.034295446658416
$ws.Cells.Item(22, 10).Value = 1.030390407367261
$ws.Cells.Item(22, 11).Value = 1.039432124313043
$ws.Cells.Item(22, 12).Value = 1.048248117019939
$ws.Cells.Item(22, 13).Value = 1.052701720882566
$ws.Cells.Item(22, 14).Value = 1.014140772222875

# Row 23 -- r="23"
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023777304519886
$ws.Cells.Item(23, 4).Value = 1.035964525654786
$ws.Cells.Item(23, 5).Value = 1.044887369309668
$ws.Cells.Item(23, 6).Value = 1.04934387022518
$ws.Cells.Item(23, 9).Value = 1.034346681103991
$ws.Cells.Item(23, 10).Value = 1.030568742235083
$ws.Cells.Item(23, 11).Value = 1.039605490943559
$ws.Cells.Item(23, 12).Value = 1.048494843870407
$ws.Cells.Item(23, 13).Value = 1.052934853807042
$ws.Cells.Item(23, 14).Value = 1.014200511282856

# Row 24 -- r="24"
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025006531621505
$ws.Cells.Item(24, 4).Value = 1.036920715976421
$ws.Cells.Item(24, 5).Value = 1.04613268207785
$ws.Cells.Item(24, 6).Value = 1.0505333298428
$ws.Cells.Item(24, 9).Value = 1.034545894723225
$ws.Cells.Item(24, 10).Value = 1.031270944422411
$ws.Cells.Item(24, 11).Value = 1.040287567559785
$ws.Cells.Item(24, 12).Value = 1.049467760594879
$ws.Cells.Item(24, 13).Value = 1.05385344624011
$ws.Cells.Item(24, 14).Value = 1.014435648368092

# Row 25 -- r="25"
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026436226334202
$ws.Cells.Item(25, 4).Value = 1.038032562321533
$ws.Cells.Item(25, 5).Value = 1.047584028122362
$ws.Cells.Item(25, 6).Value = 1.051918232210316
$ws.Cells.Item(25, 9).Value = 1.034771884265255
$ws.Cells.Item(25, 10).Value = 1.03208600071173
$ws.Cells.Item(25, 11).Value = 1.041078091786658
$ws.Cells.Item(25, 12).Value = 1.050599978318465
$ws.Cells.Item(25, 13).Value = 1.054920953388501
$ws.Cells.Item(25, 14).Value = 1.014708390798629
